# Updated fertility estimates (columnsFertilityF1a / columnsFertilityF1b columns)
# in the ColumnsNumberParameters sheet: the values are entered as text
# (quote-prefixed) rather than numbers, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")
[void]$ws.Activate()

# Row 27 -> columnsFertilityF1b = 28 (entered first so it lands earlier in
# the shared-string table, matching the authored file)
$ws.Range("B27").Value = "'28"

# Row 26 -> columnsFertilityF1a = 5
$ws.Range("B26").Value = "'5"

# Leave the selection on B27, as in the saved workbook.
[void]$ws.Range("B27").Select()
